# The document contains a single table whose "Trigger" and "erledigt"
# columns hold "Ja"/"ja"/"Nein"/"nein" values for several rows. The edit
# removes the runs holding those values (leaving the paragraphs empty),
# in document order. We locate each occurrence with Find (case sensitive,
# whole word) and remove it by replacing with an empty string, one match
# at a time (wdReplaceOne = 1), so earlier removals don't disturb the
# positions of the remaining ones.

$d = $word.ActiveDocument

$targets = @("Ja", "Ja", "Ja", "Ja", "Ja", "Ja", "ja", "Ja", "Ja", "Ja", "ja", "Nein", "Ja", "ja", "nein")

foreach ($word_text in $targets) {
    $d.Content.Find.Execute($word_text, $true, $false, $false, $false, $false, $true, 1, $false, "", 1)
}
